$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.051.97"
$ws.Range("E2").Value = "  +1.74%  "
$ws.Range("D3").Value = "3.617.36"
$ws.Range("E3").Value = "  +4.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.03"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "657.29"
$ws.Range("E6").Value = "  +4.98%  "
$ws.Range("E7").Value = "  +1.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.399"
$ws.Range("E8").Value = "  +1.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.994"
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("D11").Value = "3.615.59"
$ws.Range("E11").Value = "  +4.40%  "
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.02"
$ws.Range("E13").Value = "  -3.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.42"
$ws.Range("E14").Value = "  +2.95%  "
$ws.Range("D15").Value = "4.304.04"
$ws.Range("E15").Value = "  +4.50%  "
$ws.Range("D16").Value = "94.916.04"
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000251"
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("D18").Value = "3.598.13"
$ws.Range("E18").Value = "  +3.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.91"
$ws.Range("E19").Value = "  -4.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.85"
$ws.Range("E20").Value = "  +8.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.95"
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.53"
$ws.Range("E22").Value = "  +4.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.478"
$ws.Range("E23").Value = "  -4.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "501.64"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000195"
$ws.Range("E25").Value = "  +6.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.57"
$ws.Range("E26").Value = "  -3.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "94.47"
$ws.Range("E27").Value = "  +2.80%  "
$ws.Range("D28").Value = "3.811.10"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "12.44"
$ws.Range("E29").Value = "  +2.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.10"
$ws.Range("E30").Value = "  +9.33%  "
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.21"
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.138"
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "32.22"
$ws.Range("E35").Value = "  +9.21%  "
$ws.Range("E36").Value = "  -2.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.556"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "562.88"
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.01"
$ws.Range("E39").Value = "  +6.33%  "
$ws.Range("E40").Value = "  +1.80%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  +0.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.912"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "36.17"
$ws.Range("E44").Value = "  +45.15%  "
$ws.Range("E45").Value = "  +1.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.68"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.63"
$ws.Range("E47").Value = "  +2.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.24"
$ws.Range("E48").Value = "  +5.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0411"
$ws.Range("E49").Value = "  -3.02%  "
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.42"
$ws.Range("E51").Value = "  +0.51%  "
